# Update the "want-to-go" counts (column F) and "lowest price" (column G)
# for several con/event rows on the "展览" (Exhibitions) and "全部类型"
# (All types) sheets, reflecting freshly scraped numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("G2").Value = 35
$wsExpo.Range("F3").Value = 326
$wsExpo.Range("F4").Value = 417
$wsExpo.Range("F5").Value = 1707
$wsExpo.Range("F7").Value = 2160
$wsExpo.Range("F11").Value = 4826
$wsExpo.Range("F21").Value = 3755
$wsExpo.Range("F33").Value = 22
$wsExpo.Range("F34").Value = 873
$wsExpo.Range("F35").Value = 2384

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("G2").Value = 35
$wsAll.Range("F3").Value = 326
$wsAll.Range("F4").Value = 417
$wsAll.Range("F5").Value = 1707
$wsAll.Range("F7").Value = 2160
$wsAll.Range("F11").Value = 4826
$wsAll.Range("F21").Value = 3755
$wsAll.Range("F34").Value = 22
$wsAll.Range("F35").Value = 873
$wsAll.Range("F36").Value = 2384
